$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("U1").Value = "SH Code"

$ws.Range("U2").Value = 210610
$ws.Range("U3").Value = 210610
$ws.Range("U4").Value = 210610

$ws.Range("U4").Select()
